# Mise à jour du fichier via Shiny
#
# 1) Break the now-unused external link to "Maquette CODITRIM_test.xlsx"
#    (this removes xl/externalLinks/externalLink2.xml, renumbers the
#    remaining external links, and drops the corresponding
#    <externalReference> entry from the workbook).
# 2) Fix up the [N] external-book indices inside the surviving hidden
#    defined names so they keep pointing at the same physical external
#    workbook after the renumbering caused by step 1.
# 3) Replace the computed/linked values on "pro" and "VA" with the new
#    hard values pulled in by the refresh, and update each sheet's
#    selection to reflect where the user left the cursor (B2:B26).

$wb = $excel.ActiveWorkbook

# --- 1) Break the stale external link ------------------------------------
$wb.BreakLink("file:///C:\Users\HP\Documents\Maquette%20CODITRIM_test.xlsx", 1)

# --- 2) Re-point the surviving [N] external references --------------------
$wb.Names.Item("__123Graph_D").RefersTo = "=[2]E!#REF!"
$wb.Names.Item("__123Graph_E").RefersTo = "=[2]E!#REF!"
$wb.Names.Item("graphe").RefersTo = "=[2]E!#REF!"
$wb.Names.Item("_123graph_b").RefersTo = "=[3]A!#REF!"
$wb.Names.Item("pol").RefersTo = "=[3]A!#REF!"
$wb.Names.Item("_FilterDatabase").RefersTo = "=[4]C!`$P`$428:`$T`$428"
$wb.Names.Item("_Regression_Out").RefersTo = "=[4]C!`$AK`$18:`$AK`$18"
$wb.Names.Item("_Regression_X").RefersTo = "=[4]C!`$AK`$11:`$AU`$11"
$wb.Names.Item("_Regression_Y").RefersTo = "=[4]C!`$AK`$10:`$AU`$10"
$wb.Names.Item("ACwvu.Print.").RefersTo = "=[5]Med!#REF!"
$wb.Names.Item("Swvu.Print.").RefersTo = "=[5]Med!#REF!"
$wb.Names.Item("CRISa").RefersTo = "=[6]A!#REF!"

# --- 3) Refresh the data values --------------------------------------------
$newValues = @(689569,686496,726689,597520,622389,686246,819599,870647,990331,953418,985082,1083865,1118119,1187401,1155632,929833,1403276,1590431,1834812,2000000,2080000,2287000,2411000,2714431,2463842)

$wsPro = $wb.Worksheets.Item("pro")
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $wsPro.Cells.Item($i + 2, 2).Value = $newValues[$i]
}
$wsPro.Range("B2:B26").Select()

$wsVA = $wb.Worksheets.Item("VA")
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $wsVA.Cells.Item($i + 2, 2).Value = $newValues[$i]
}
$wsVA.Range("B2:B26").Select()

# Restore "pro" as the active sheet/selection (it was tabSelected before the
# edit, and selecting on "VA" above necessarily activated that sheet).
$wsPro.Activate()
$wsPro.Range("B2:B26").Select()
